# Add "zmin" / "zmax" columns (J, K) to Sheet1, matching the missing
# columns that several other expdata files already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new column headers, centered like the rest of row 1.
$ws.Range("J1").Value = "zmin"
$ws.Range("K1").Value = "zmax"
$ws.Range("J1:K1").HorizontalAlignment = -4108   # xlCenter

# Data rows 2-13: constant zmin = 0.2, zmax = 1 for every data point.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.2
    $ws.Cells.Item($r, 11).Value = 1
}

# Move the active selection to L11, matching the saved workbook state.
$ws.Range("L11").Select() | Out-Null
